$wb = $excel.ActiveWorkbook
$wsOut = $wb.Worksheets.Item("2aOutput")
$ws4 = $wb.Worksheets.Item("Sheet4")

$ws4.Range("A1").Value = "NumDocs"
$ws4.Range("B1").Value = "busyness"

for ($i = 2; $i -le 21; $i++) {
    $ws4.Cells.Item($i, 1).Value = $wsOut.Cells.Item($i, 1).Value2
    $ws4.Cells.Item($i, 2).Value = $wsOut.Cells.Item($i, 10).Value2
}

$co = $ws4.ChartObjects().Add(137.875, 50.625, 433.0625, 216)
$chart = $co.Chart
$chart.ChartType = 72
$ser = $chart.SeriesCollection().NewSeries()
$ser.Formula = "=SERIES(Sheet4!`$B`$1,Sheet4!`$A`$2:`$A`$21,Sheet4!`$B`$2:`$B`$21,1)"
$ser.Smooth = $true

$ws4.Activate()

Write-Host "done"
